$d = $word.ActiveDocument

# --- Change 1: "- Tweaked player jumping speed, physics etc to improve overall feel"
# is split into three runs, with a proofErr spellStart/spellEnd pair wrapping "etc"
# (a cosmetic spell-checker marker; the visible text is unchanged).
$p7 = $d.Paragraphs.Item(7).Range
if ($p7.Text -notmatch "Tweaked player jumping speed") {
    throw "Paragraph 7 is not the 'Tweaked player jumping speed' bullet; document structure differs from what was expected."
}
$xmlTweak = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">- Tweaked player jumping speed, physics </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>etc</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> to improve overall feel</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$p7.InsertXML($xmlTweak)

# --- Change 2: rework the "Small Bugs" bullet list.
# Original order (paragraphs 32-36 of the "Small Bugs" list):
#   32 Enemy collision detection...
#   33 The player can outrun the camera...
#   34 We left out visual cues...
#   35 Animation for kick...
#   36 When fighting the boss... (+ _GoBack bookmark)     <- last paragraph in the doc body
#
# Target order (3 bullets only):
#   A  The player can outrun the camera...                          (unchanged text)
#   B  (bookmark moved here) Animation for spin/kick... (spin+kick proofErr'd)
#   C  When fighting the boss... (NullReferenceException / functions proofErr'd, bookmark removed)

if ($d.Paragraphs.Item(32).Range.Text -notmatch "Enemy collision detection") {
    throw "Paragraph 32 is not the 'Enemy collision detection' bullet; document structure differs from what was expected."
}
if ($d.Paragraphs.Item(33).Range.Text -notmatch "The player can outrun the camera") {
    throw "Paragraph 33 is not the 'The player can outrun the camera' bullet; document structure differs from what was expected."
}
if ($d.Paragraphs.Item(34).Range.Text -notmatch "We left out visual cues") {
    throw "Paragraph 34 is not the 'We left out visual cues' bullet; document structure differs from what was expected."
}
if ($d.Paragraphs.Item(35).Range.Text -notmatch "Animation for kick") {
    throw "Paragraph 35 is not the 'Animation for kick' bullet; document structure differs from what was expected."
}
if ($d.Paragraphs.Item(36).Range.Text -notmatch "When fighting the boss") {
    throw "Paragraph 36 is not the 'When fighting the boss' bullet; document structure differs from what was expected."
}

# Remove the "Enemy collision detection..." bullet entirely.
$d.Paragraphs.Item(32).Range.Delete()

# Remove the "We left out visual cues..." bullet entirely (it is now at index 33).
$d.Paragraphs.Item(33).Range.Delete()

# Now exactly 3 bullets remain, at indices 32, 33, 34:
#   32 The player can outrun the camera...     -> stays as-is (matches target A)
#   33 Animation for kick...                   -> replace with bookmark + spin/kick text (target B)
#   34 When fighting the boss...(+bookmark)     -> replace with proofErr'd text, drop bookmark (target C)

$pB = $d.Paragraphs.Item(33).Range
$xmlB = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="4"/>
              </w:numPr>
            </w:pPr>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r><w:t xml:space="preserve">Animation for </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>spin</w:t></w:r>
            <w:r><w:t>kick</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> doesn&#8217;t start until player hits the ground but to try and fix it caused many other problems so it&#8217;s still in there.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$pB.InsertXML($xmlB)

# The third bullet (currently paragraph 34) is the LAST paragraph in the document body.
# Replacing the last body paragraph's content via InsertXML leaves a stray empty paragraph
# behind (the trailing paragraph mark does not get consumed), so instead: insert a new
# paragraph before it (so it is no longer last), fill that new one via InsertXML, then
# delete the original.
$pOldLast = $d.Paragraphs.Item(34).Range
$pOldLast.InsertParagraphBefore()

$pC = $d.Paragraphs.Item(34).Range
$xmlC = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="4"/>
              </w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">When fighting the boss, there is a </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>NullReferenceException</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> however there is no reason for it because the game still </w:t></w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r><w:t>functions</w:t></w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r><w:t xml:space="preserve"> as it should.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@
$pC.InsertXML($xmlC)

# Delete the now-redundant original last paragraph (it still has the old text + bookmark).
$d.Paragraphs.Item(35).Range.Delete()

Write-Host "Edit complete."
